$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 225
$ws.Cells.Item(12, 9).Value = 225
$ws.Cells.Item(12, 11).Value = 225
$ws.Cells.Item(12, 13).Value = -55
$ws.Cells.Item(17, 8).Value = 1309.125
$ws.Cells.Item(17, 10).Value = 1333.2858
$ws.Cells.Item(17, 12).Value = 3999.8574
$ws.Cells.Item(17, 14).Value = -4335.857400000001
$ws.Cells.Item(46, 8).Value = 169833.33
$ws.Cells.Item(46, 10).Value = 203200
$ws.Cells.Item(46, 12).Value = 609600
$ws.Cells.Item(46, 14).Value = -609838
$ws.Cells.Item(60, 8).Value = 169833.33
$ws.Cells.Item(60, 10).Value = 203200
$ws.Cells.Item(60, 12).Value = 609600
$ws.Cells.Item(60, 14).Value = -610568
$ws.Cells.Item(125, 8).Value = 6299
$ws.Cells.Item(125, 10).Value = 4060
$ws.Cells.Item(125, 12).Value = 36540
$ws.Cells.Item(125, 14).Value = -41460
$ws.Cells.Item(133, 8).Value = 99333
$ws.Cells.Item(133, 10).Value = 99333
$ws.Cells.Item(133, 12).Value = 99333
$ws.Cells.Item(133, 14).Value = -109453
$ws.Cells.Item(135, 8).Value = 4335.4062
$ws.Cells.Item(135, 9).Value = 4507.4443
$ws.Cells.Item(135, 11).Value = 40566.9987
$ws.Cells.Item(135, 13).Value = -38031.9987

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7979.041
$ws.Cells.Item(32, 9).Value = 7663.8955
$ws.Cells.Item(32, 11).Value = 7663.8955
$ws.Cells.Item(32, 13).Value = -7376.8955
$ws.Cells.Item(61, 8).Value = 7790.094
$ws.Cells.Item(61, 9).Value = 6929.0454
$ws.Cells.Item(61, 11).Value = 6929.0454
$ws.Cells.Item(61, 13).Value = -6717.0454
$ws.Cells.Item(74, 8).Value = 2017.8928
$ws.Cells.Item(74, 10).Value = 2073.1155
$ws.Cells.Item(74, 12).Value = 2073.1155
$ws.Cells.Item(74, 14).Value = -3821.1155
$ws.Cells.Item(77, 8).Value = 2017.8928
$ws.Cells.Item(77, 10).Value = 2073.1155
$ws.Cells.Item(77, 12).Value = 10365.5775
$ws.Cells.Item(77, 14).Value = -19101.5775
$ws.Cells.Item(136, 8).Value = 7790.094
$ws.Cells.Item(136, 9).Value = 6929.0454
$ws.Cells.Item(136, 11).Value = 20787.1362
$ws.Cells.Item(136, 13).Value = -18237.1362

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 9384.6
$ws.Cells.Item(20, 10).Value = 11099.8
$ws.Cells.Item(20, 12).Value = 11099.8
$ws.Cells.Item(20, 14).Value = -11593.8
$ws.Cells.Item(54, 8).Value = 7750
$ws.Cells.Item(54, 9).Value = 7750
$ws.Cells.Item(54, 11).Value = 7750
$ws.Cells.Item(54, 13).Value = -7266
$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 14).Value = ""
$ws.Cells.Item(81, 8).Value = 106796.5
$ws.Cells.Item(81, 10).Value = 106796.5
$ws.Cells.Item(81, 12).Value = 106796.5
$ws.Cells.Item(81, 14).Value = -108918.5
$ws.Cells.Item(84, 8).Value = 106796.5
$ws.Cells.Item(84, 10).Value = 106796.5
$ws.Cells.Item(84, 12).Value = 320389.5
$ws.Cells.Item(84, 14).Value = -330997.5
$ws.Cells.Item(102, 8).Value = 28274.5
$ws.Cells.Item(102, 9).Value = 7550
$ws.Cells.Item(102, 10).Value = 48999
$ws.Cells.Item(102, 11).Value = 7550
$ws.Cells.Item(102, 12).Value = 48999
$ws.Cells.Item(102, 13).Value = -4305
$ws.Cells.Item(102, 14).Value = -55489
$ws.Cells.Item(105, 8).Value = 6050.643
$ws.Cells.Item(105, 9).Value = 10008.917
$ws.Cells.Item(105, 10).Value = 3081.9375
$ws.Cells.Item(105, 11).Value = 10008.917
$ws.Cells.Item(105, 12).Value = 3081.9375
$ws.Cells.Item(105, 13).Value = -8261.916999999999
$ws.Cells.Item(105, 14).Value = -6575.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 9504.762000000001
$ws.Cells.Item(62, 9).Value = 7900.643
$ws.Cells.Item(62, 11).Value = 7900.643
$ws.Cells.Item(62, 13).Value = -7276.643
$ws.Cells.Item(65, 8).Value = 9504.762000000001
$ws.Cells.Item(65, 9).Value = 7900.643
$ws.Cells.Item(65, 11).Value = 39503.215
$ws.Cells.Item(65, 13).Value = -36383.215
$ws.Cells.Item(99, 8).Value = 6217147.5
$ws.Cells.Item(99, 10).Value = 4354.579
$ws.Cells.Item(99, 12).Value = 4354.579
$ws.Cells.Item(99, 14).Value = -7350.579
$ws.Cells.Item(109, 8).Value = 59994
$ws.Cells.Item(109, 10).Value = 59994
$ws.Cells.Item(109, 12).Value = 59994
$ws.Cells.Item(109, 14).Value = -62074
$ws.Cells.Item(114, 8).Value = 52554
$ws.Cells.Item(114, 10).Value = 52554
$ws.Cells.Item(114, 12).Value = 52554
$ws.Cells.Item(114, 14).Value = -61232
$ws.Cells.Item(126, 8).Value = 6217147.5
$ws.Cells.Item(126, 10).Value = 4354.579
$ws.Cells.Item(126, 12).Value = 13063.737
$ws.Cells.Item(126, 14).Value = -18003.737
$ws.Cells.Item(141, 8).Value = 201894.81
$ws.Cells.Item(141, 10).Value = 216208.34
$ws.Cells.Item(141, 12).Value = 216208.34
$ws.Cells.Item(141, 14).Value = -226568.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(59, 8).Value = 2551
$ws.Cells.Item(59, 10).Value = 2000
$ws.Cells.Item(59, 12).Value = 6000
$ws.Cells.Item(59, 14).Value = -7080
$ws.Cells.Item(62, 8).Value = 2000
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 2000
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 6000
$ws.Cells.Item(62, 13).Value = ""
$ws.Cells.Item(62, 14).Value = -7372
$ws.Cells.Item(65, 8).Value = 2000
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 2000
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 18000
$ws.Cells.Item(65, 13).Value = ""
$ws.Cells.Item(65, 14).Value = -24864

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8493.414000000001
$ws.Cells.Item(70, 9).Value = 6837.6
$ws.Cells.Item(70, 10).Value = 10267.5
$ws.Cells.Item(70, 11).Value = 6837.6
$ws.Cells.Item(70, 12).Value = 10267.5
$ws.Cells.Item(70, 13).Value = -6567.6
$ws.Cells.Item(70, 14).Value = -10807.5
$ws.Cells.Item(73, 8).Value = 8493.414000000001
$ws.Cells.Item(73, 9).Value = 6837.6
$ws.Cells.Item(73, 10).Value = 10267.5
$ws.Cells.Item(73, 11).Value = 6837.6
$ws.Cells.Item(73, 12).Value = 10267.5
$ws.Cells.Item(73, 13).Value = -5901.6
$ws.Cells.Item(73, 14).Value = -12139.5
$ws.Cells.Item(126, 8).Value = 17890.285
$ws.Cells.Item(126, 9).Value = 20568
$ws.Cells.Item(126, 11).Value = 61704
$ws.Cells.Item(126, 13).Value = -59234

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 16045.84
$ws.Cells.Item(7, 9).Value = 18590.475
$ws.Cells.Item(7, 11).Value = 18590.475
$ws.Cells.Item(7, 13).Value = -18478.475
$ws.Cells.Item(22, 8).Value = 806.95
$ws.Cells.Item(22, 9).Value = 676.75
$ws.Cells.Item(22, 10).Value = 893.75
$ws.Cells.Item(22, 11).Value = 676.75
$ws.Cells.Item(22, 12).Value = 893.75
$ws.Cells.Item(22, 13).Value = -381.75
$ws.Cells.Item(22, 14).Value = -1483.75
$ws.Cells.Item(27, 8).Value = 806.95
$ws.Cells.Item(27, 9).Value = 676.75
$ws.Cells.Item(27, 10).Value = 893.75
$ws.Cells.Item(27, 11).Value = 676.75
$ws.Cells.Item(27, 12).Value = 893.75
$ws.Cells.Item(27, 13).Value = -569.75
$ws.Cells.Item(27, 14).Value = -1107.75
$ws.Cells.Item(46, 8).Value = 1338.8636
$ws.Cells.Item(46, 9).Value = 772.75
$ws.Cells.Item(46, 11).Value = 772.75
$ws.Cells.Item(46, 13).Value = -584.75
$ws.Cells.Item(63, 8).Value = 68000
$ws.Cells.Item(63, 10).Value = 68000
$ws.Cells.Item(63, 12).Value = 68000
$ws.Cells.Item(63, 14).Value = -69498
$ws.Cells.Item(66, 8).Value = 68000
$ws.Cells.Item(66, 10).Value = 68000
$ws.Cells.Item(66, 12).Value = 204000
$ws.Cells.Item(66, 14).Value = -211488
$ws.Cells.Item(82, 8).Value = 1639.8334
$ws.Cells.Item(82, 10).Value = 2339.8
$ws.Cells.Item(82, 12).Value = 2339.8
$ws.Cells.Item(82, 14).Value = -3061.8
$ws.Cells.Item(85, 8).Value = 1639.8334
$ws.Cells.Item(85, 10).Value = 2339.8
$ws.Cells.Item(85, 12).Value = 2339.8
$ws.Cells.Item(85, 14).Value = -4835.8
$ws.Cells.Item(123, 8).Value = 107489
$ws.Cells.Item(123, 10).Value = 107489
$ws.Cells.Item(123, 12).Value = 107489
$ws.Cells.Item(123, 14).Value = -117289
$ws.Cells.Item(126, 8).Value = 16045.84
$ws.Cells.Item(126, 9).Value = 18590.475
$ws.Cells.Item(126, 11).Value = 55771.425
$ws.Cells.Item(126, 13).Value = -53301.425
$ws.Cells.Item(129, 8).Value = 77777
$ws.Cells.Item(129, 10).Value = 77777
$ws.Cells.Item(129, 12).Value = 77777
$ws.Cells.Item(129, 14).Value = -87777
$ws.Cells.Item(136, 8).Value = 6994.4375
$ws.Cells.Item(136, 10).Value = 6916.385
$ws.Cells.Item(136, 12).Value = 20749.155
$ws.Cells.Item(136, 14).Value = -25849.155

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1371.9487
$ws.Cells.Item(113, 9).Value = 851.7273
$ws.Cells.Item(113, 11).Value = 2555.1819
$ws.Cells.Item(113, 13).Value = -385.1819
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).Value = ""
$ws.Cells.Item(126, 8).Value = 35007.46
$ws.Cells.Item(126, 9).Value = 82378.39999999999
$ws.Cells.Item(126, 10).Value = 5400.625
$ws.Cells.Item(126, 11).Value = 247135.2
$ws.Cells.Item(126, 12).Value = 16201.875
$ws.Cells.Item(126, 13).Value = -244665.2
$ws.Cells.Item(126, 14).Value = -21141.875
$ws.Cells.Item(136, 8).Value = 329687.97
$ws.Cells.Item(136, 9).Value = 336724.66
$ws.Cells.Item(136, 11).Value = 1010173.98
$ws.Cells.Item(136, 13).Value = -1007623.98
